$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark from the paragraph that ends
#    with "...risk factor of 33 percent."
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 2. After that paragraph, insert: an empty paragraph, a new paragraph with
#    the new answer text (tab + text), and another empty paragraph.
$pRisk = $d.Paragraphs.Item(17)
$pRisk.Range.InsertParagraphAfter()

$pBlank1 = $d.Paragraphs.Item(18)
$pBlank1.Range.InsertParagraphAfter()

$pNewText = $d.Paragraphs.Item(19)
$pNewText.Range.InsertParagraphAfter()

$pNewText = $d.Paragraphs.Item(19)
$pNewText.Range.InsertAfter("`tThere are three individual trips needed with one event that presents danger. I believe the only possible solution is leaving the cat and parrot alone during these trips. The seed has no defense against the parrot. While the cat and parrot can fight it out while left alone. If they are left alone for say 30 minutes there is a chance both parties would survive the fight or maybe complete peace among these two would prevail.")

# 3. Move the "_GoBack" bookmark to the (pre-existing) empty paragraph that
#    follows the "Socks in the dark" heading.
$pAfterSocks = $d.Paragraphs.Item(22)
$d.Bookmarks.Add("_GoBack", $pAfterSocks.Range)
